$d = $word.ActiveDocument

$replacements = @(
    @{old = "55×16=880"; new = "42×37=1554"},
    @{old = "16×14=224"; new = "86×34=2924"},
    @{old = "19×12=228"; new = "12×39=468"},
    @{old = "55×77=4235"; new = "64×28=1792"},
    @{old = "31×42=1302"; new = "41×80=3280"},
    @{old = "73×25=1825"; new = "14×16=224"},
    @{old = "42×49=2058"; new = "79×57=4503"},
    @{old = "37×67=2479"; new = "25×55=1375"},
    @{old = "86×60=5160"; new = "72×60=4320"},
    @{old = "25×66=1650"; new = "42×31=1302"},
    @{old = "80×19=1520"; new = "34×21=714"},
    @{old = "19×77=1463"; new = "67×96=6432"},
    @{old = "21×65=1365"; new = "79×88=6952"},
    @{old = "30×60=1800"; new = "99×59=5841"},
    @{old = "71×95=6745"; new = "85×13=1105"},
    @{old = "14×65=910"; new = "50×57=2850"},
    @{old = "66×63=4158"; new = "64×22=1408"},
    @{old = "85×79=6715"; new = "81×91=7371"},
    @{old = "99×55=5445"; new = "71×67=4757"},
    @{old = "49×63=3087"; new = "53×38=2014"},
    @{old = "77×72=5544"; new = "33×53=1749"},
    @{old = "47×70=3290"; new = "69×67=4623"},
    @{old = "89×18=1602"; new = "95×47=4465"},
    @{old = "27×16=432"; new = "11×59=649"},
    @{old = "43×30=1290"; new = "64×51=3264"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
